# Pio's ERA operativo Abril-Diciembre 2025
# Update "Antiguedad" (C) and "Carga Teorica" (F) figures for all medicos,
# plus a handful of monthly load ratios (O2, P10, R13), reflecting the
# refreshed reference date used to compute seniority/workload.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Gomez
$ws.Cells.Item(2, 3).Value = 10.92328767123288
$ws.Cells.Item(2, 6).Value = 6.224379325864565
$ws.Cells.Item(2, 15).Value = 0.25

# Row 3 - Bravo
$ws.Cells.Item(3, 3).Value = 10.92328767123288
$ws.Cells.Item(3, 6).Value = 6.224379325864565

# Row 4 - Iñiguez
$ws.Cells.Item(4, 3).Value = 10.41917808219178
$ws.Cells.Item(4, 6).Value = 6.280334269970441

# Row 5 - Breinbauer
$ws.Cells.Item(5, 3).Value = 10.41917808219178
$ws.Cells.Item(5, 6).Value = 6.280334269970441

# Row 6 - Arredondo
$ws.Cells.Item(6, 3).Value = 9.838356164383562
$ws.Cells.Item(6, 6).Value = 6.344804096875039

# Row 7 - Carrasco
$ws.Cells.Item(7, 3).Value = 9.838356164383562
$ws.Cells.Item(7, 6).Value = 6.344804096875039

# Row 8 - Culaciati
$ws.Cells.Item(8, 3).Value = 9.838356164383562
$ws.Cells.Item(8, 6).Value = 6.344804096875039

# Row 9 - Contreras
$ws.Cells.Item(9, 3).Value = 7.416438356164384
$ws.Cells.Item(9, 6).Value = 6.613631110948924

# Row 10 - Cisternas
$ws.Cells.Item(10, 3).Value = 6.583561643835616
$ws.Cells.Item(10, 6).Value = 6.706078409906459
$ws.Cells.Item(10, 16).Value = 0.25

# Row 11 - Pio
$ws.Cells.Item(11, 3).Value = 6.383561643835616
$ws.Cells.Item(11, 6).Value = 6.728277925774551

# Row 12 - Alvo
$ws.Cells.Item(12, 3).Value = 4.583561643835616
$ws.Cells.Item(12, 6).Value = 6.928073568587381

# Row 13 - Boettiger
$ws.Cells.Item(13, 3).Value = 2.747945205479452
$ws.Cells.Item(13, 6).Value = 7.131822549842476
$ws.Cells.Item(13, 18).Value = 0.5

# Row 14 - Loch
$ws.Cells.Item(14, 3).Value = 2
$ws.Cells.Item(14, 6).Value = 7.214842657129998

# Row 15 - Rubio
$ws.Cells.Item(15, 3).Value = 2
$ws.Cells.Item(15, 6).Value = 7.214842657129998

# Row 16 - Recluta1
$ws.Cells.Item(16, 3).Value = 0.1643835616438356
$ws.Cells.Item(16, 6).Value = 7.418591638385093
